# Update "想去人数" (number of people interested) values on the
# "展览" and "全部类型" worksheets.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 75
$ws1.Range("F6").Value = 5586
$ws1.Range("F8").Value = 6461
$ws1.Range("F10").Value = 17
$ws1.Range("F11").Value = 1415
$ws1.Range("F13").Value = 41
$ws1.Range("F14").Value = 112

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 75
$ws4.Range("F7").Value = 5586
$ws4.Range("F9").Value = 6461
$ws4.Range("F11").Value = 17
$ws4.Range("F12").Value = 1415
$ws4.Range("F14").Value = 41
$ws4.Range("F15").Value = 112
